$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCDTtiNTY")

# Row 1: remove the "Share that is New (dimensionless)" header cell/note and
# restore the default (unwrapped) row height.
$ws.Range("A1").Clear()
$ws.Rows.Item(1).AutoFit()

# Row 2 - LDVs: updated passenger/freight shares
$ws.Range("B2").Value = 0.0815
$ws.Range("C2").Value = 0.083

# Row 3 - HDVs: updated passenger/freight shares
$ws.Range("B3").Value = 0.065
$ws.Range("C3").Value = 0.066

# Row 4 - aircraft: freight share updated (passenger unchanged)
$ws.Range("C4").Value = 0.0425

# Row 5 - rail: passenger share updated (freight unchanged)
$ws.Range("B5").Value = 0.043

# Row 6 - ships: unchanged, left as-is

# Row 7 - motorbikes: passenger share updated, freight becomes a literal
# value instead of a formula, and the explanatory note in D7 is removed.
$ws.Range("B7").Value = 0.14
$ws.Range("C7").Value = 0.0825
$ws.Range("D7").Clear()

$ws.Range("C8").Select() | Out-Null
